$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 91.34108066666666
$ws.Range("H2").Value = 274.023242
$ws.Range("I2").Value = 0.2190334467302001
$ws.Range("J2").Value = 0.2190334467302
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2615913333333333
$ws.Range("N2").Value = 0.784774
$ws.Range("O2").Value = 0.08239613548481725
$ws.Range("P2").Value = 0.08239613548481727
$ws.Range("Q2").Value = 23.89403507970089
$ws.Range("R2").Value = 215.046315717308
$ws.Range("S2").Value = 0.01804750955248807
$ws.Range("T2").Value = 0.01804750955248807

# Row 3
$ws.Range("G3").Value = 91.34108066666666
$ws.Range("H3").Value = 274.023242
$ws.Range("I3").Value = 0.2190334467302001
$ws.Range("J3").Value = 0.2190334467302
$ws.Range("N3").Value = 5.233242000000001
$ws.Range("O3").Value = 0.5494561706387266
$ws.Range("P3").Value = 0.5494561706387268
$ws.Range("Q3").Value = 159.3366598900627
$ws.Range("R3").Value = 1434.029939010564
$ws.Range("S3").Value = 0.1203492788821773
$ws.Range("T3").Value = 0.1203492788821773

# Row 4
$ws.Range("G4").Value = 91.34108066666666
$ws.Range("H4").Value = 274.023242
$ws.Range("I4").Value = 0.2190334467302001
$ws.Range("J4").Value = 0.2190334467302
$ws.Range("M4").Value = 1.168795666666667
$ws.Range("N4").Value = 3.506387
$ws.Range("O4").Value = 0.3681476938764561
$ws.Range("P4").Value = 0.3681476938764561
$ws.Range("Q4").Value = 106.7590592718504
$ws.Range("R4").Value = 960.831533446654
$ws.Range("S4").Value = 0.08063665829553474
$ws.Range("T4").Value = 0.08063665829553474

# Row 5
$ws.Range("G5").Value = 276.4348856666666
$ws.Range("H5").Value = 829.3046569999999
$ws.Range("I5").Value = 0.6628833966285105
$ws.Range("J5").Value = 0.6628833966285105
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2615913333333333
$ws.Range("N5").Value = 0.784774
$ws.Range("O5").Value = 0.08239613548481725
$ws.Range("P5").Value = 0.08239613548481727
$ws.Range("Q5").Value = 72.31297032139088
$ws.Range("R5").Value = 650.8167328925178
$ws.Range("S5").Value = 0.0546190301592386
$ws.Range("T5").Value = 0.05461903015923861

# Row 6
$ws.Range("G6").Value = 276.4348856666666
$ws.Range("H6").Value = 829.3046569999999
$ws.Range("I6").Value = 0.6628833966285105
$ws.Range("J6").Value = 0.6628833966285105
$ws.Range("N6").Value = 5.233242000000001
$ws.Range("O6").Value = 0.5494561706387266
$ws.Range("P6").Value = 0.5494561706387268
$ws.Range("Q6").Value = 482.2168846453326
$ws.Range("R6").Value = 4339.951961807994
$ws.Range("S6").Value = 0.3642253726914936
$ws.Range("T6").Value = 0.3642253726914936

# Row 7
$ws.Range("G7").Value = 276.4348856666666
$ws.Range("H7").Value = 829.3046569999999
$ws.Range("I7").Value = 0.6628833966285105
$ws.Range("J7").Value = 0.6628833966285105
$ws.Range("M7").Value = 1.168795666666667
$ws.Range("N7").Value = 3.506387
$ws.Range("O7").Value = 0.3681476938764561
$ws.Range("P7").Value = 0.3681476938764561
$ws.Range("Q7").Value = 323.0958964826954
$ws.Range("R7").Value = 2907.863068344259
$ws.Range("S7").Value = 0.2440389937777783
$ws.Range("T7").Value = 0.2440389937777783

# Row 8
$ws.Range("G8").Value = 49.24290466666667
$ws.Range("H8").Value = 147.728714
$ws.Range("I8").Value = 0.1180831566412894
$ws.Range("J8").Value = 0.1180831566412894
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2615913333333333
$ws.Range("N8").Value = 0.784774
$ws.Range("O8").Value = 0.08239613548481725
$ws.Range("P8").Value = 0.08239613548481727
$ws.Range("Q8").Value = 12.88151708895956
$ws.Range("R8").Value = 115.933653800636
$ws.Range("S8").Value = 0.009729595773090583
$ws.Range("T8").Value = 0.009729595773090582

# Row 9
$ws.Range("G9").Value = 49.24290466666667
$ws.Range("H9").Value = 147.728714
$ws.Range("I9").Value = 0.1180831566412894
$ws.Range("J9").Value = 0.1180831566412894
$ws.Range("N9").Value = 5.233242000000001
$ws.Range("O9").Value = 0.5494561706387266
$ws.Range("P9").Value = 0.5494561706387268
$ws.Range("Q9").Value = 85.90001230119867
$ws.Range("R9").Value = 773.100110710788
$ws.Range("S9").Value = 0.06488151906505582
$ws.Range("T9").Value = 0.06488151906505582

# Row 10
$ws.Range("G10").Value = 49.24290466666667
$ws.Range("H10").Value = 147.728714
$ws.Range("I10").Value = 0.1180831566412894
$ws.Range("J10").Value = 0.1180831566412894
$ws.Range("M10").Value = 1.168795666666667
$ws.Range("N10").Value = 3.506387
$ws.Range("O10").Value = 0.3681476938764561
$ws.Range("P10").Value = 0.3681476938764561
$ws.Range("Q10").Value = 57.55489358847979
$ws.Range("R10").Value = 517.994042296318
$ws.Range("S10").Value = 0.04347204180314304
$ws.Range("T10").Value = 0.04347204180314303

